$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 7 with the "Tiebreak" round/format entry
$ws.Range("A7").Value = "Tiebreak"
$ws.Range("C7").Value = "Tiebreak"
$ws.Range("E7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("I7").Value = 0

# Update selection to match the post-edit workbook state
$ws.Application.Goto($ws.Range("J7"))
